$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Kenneth A.%Egol%NULL%2,           Sanjit R.%Konda%NULL%2,           Mackenzie L.%Bird%NULL%2,           Nicket%Dedhia%NULL%2,           Emma K.%Landes%NULL%2,           Rachel A.%Ranson%NULL%2,           Sara J.%Solasz%NULL%2,           Vinay K.%Aggarwal%NULL%2,           Joseph A.%Bosco%NULL%2,           David L.%Furgiuele%NULL%2,           Abhishek%Ganta%NULL%2,           Jason%Gould%NULL%2,           Thomas R.%Lyon%NULL%2,           Toni M.%McLaurin%NULL%2,           Nirmal C.%Tejwani%NULL%2,           Joseph D.%Zuckerman%NULL%2,           Philipp%Leucht%NULL%2]"

$ws.Range("E3").Value = "[Drake G.%LeBrun%NULL%2,           Maxwell A.%Konnaris%NULL%2,           Gregory C.%Ghahramani%NULL%2,           Ajay%Premkumar%NULL%2,           Chris J.%DeFrancesco%NULL%2,           Jordan A.%Gruskay%NULL%2,           Aleksey%Dvorzhinskiy%NULL%2,           Milan S.%Sandhu%NULL%2,           Elan M.%Goldwyn%NULL%2,           Christopher L.%Mendias%NULL%2,           William M.%Ricci%NULL%2]"

$ws.Range("E4").Value = "[Josep Maria%Muñoz Vives%NULL%1,           Montsant%Jornet-Gibert%NULL%2,           Montsant%Jornet-Gibert%NULL%0,           J.%Cámara-Cabrera%NULL%2,           J.%Cámara-Cabrera%NULL%0,           Pedro L.%Esteban%NULL%2,           Pedro L.%Esteban%NULL%0,           Laia%Brunet%NULL%2,           Laia%Brunet%NULL%0,           Luis%Delgado-Flores%NULL%2,           Luis%Delgado-Flores%NULL%0,           P.%Camacho-Carrasco%NULL%2,           P.%Camacho-Carrasco%NULL%0,           P.%Torner%NULL%2,           P.%Torner%NULL%0,           Francesc%Marcano-Fernández%NULL%2,           Francesc%Marcano-Fernández%NULL%0,           NULL%NULL%NULL%0,           NULL%NULL%NULL%0]"

$ws.Range("E5").Value = "[Zoe B.%Cheung%zoe.cheung@mountsinai.org%1,           David A.%Forsh%NULL%1]"

$ws.Range("E6").Value = "[Andrew J.%Hall%xref no email%2,    Nicholas D.%Clement%xref no email%2,    Luke%Farrow%xref no email%2,    Alasdair M. J.%MacLullich%xref no email%2,    Graham F.%Dall%xref no email%2,    Chloe E. H.%Scott%xref no email%2,    Paul J.%Jenkins%xref no email%2,    Timothy O.%White%xref no email%2,    Andrew D.%Duckworth%xref no email%2]"

$ws.Range("E7").Value = "[Babar%Kayani%xref no email%2,    Elliot%Onochie%xref no email%2,    Vijay%Patil%xref no email%2,    Fahima%Begum%xref no email%2,    Rory%Cuthbert%xref no email%2,    David%Ferguson%xref no email%2,    Jagmeet S.%Bhamra%xref no email%2,    Aadhar%Sharma%xref no email%2,    Peter%Bates%xref no email%2,    Fares S.%Haddad%xref no email%2]"
